$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new time log entry: 11/02/2023, Internship, daily operations description
# Match formatting of the existing date column (row above) before setting the
# new date value so Excel doesn't invent a new datetime-style number format.
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A18").Value = 45232

$ws.Range("B18").Value = "Internship"
$ws.Range("C18").Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

# Move the active selection to the next empty row, matching Excel's default
# behavior after data entry on the last row of a contiguous table.
$ws.Range("C19").Select()
